# Collapse the split "Suggested Retail AU: ${Retail" / " " / "AU}" runs
# in the pricing textbox on slide 1 into a single run reading
# "Suggested Retail AU: ${Retail AUD}".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TextBox 18" is shape #4 on the slide - it holds the Item#, Dims,
# FOB and Suggested Retail lines.
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

$null = $tr.Replace('Suggested Retail AU: ${Retail AU}', 'Suggested Retail AU: ${Retail AUD}', 0, 0, 0)
